$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1889763779527559
$ws.Range("C2").Value = 0.5748031496062992
$ws.Range("J2").Value = 0.01181102362204724
$ws.Range("P2").Value = 0.1574803149606299
$ws.Range("S2").Value = 0.06692913385826772
$ws.Range("B3").Value = 0.006289308176100629
$ws.Range("C3").Value = 0.08176100628930817
$ws.Range("J3").Value = 0.03144654088050314
$ws.Range("P3").Value = 0.7169811320754716
$ws.Range("S3").Value = 0.1635220125786163
$ws.Range("J4").Value = 0.06382978723404255
$ws.Range("P4").Value = 0.7021276595744681
$ws.Range("S4").Value = 0.2340425531914894
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.08290155440414508
$ws.Range("D6").Value = 0.02072538860103627
$ws.Range("F6").Value = 0.04663212435233161
$ws.Range("J6").Value = 0.2487046632124352
$ws.Range("O6").Value = 0.02072538860103627
$ws.Range("Q6").Value = 0.1709844559585492
$ws.Range("R6").Value = 0.09844559585492228
$ws.Range("S6").Value = 0.310880829015544
$ws.Range("B7").Value = 0.08936170212765958
$ws.Range("D7").Value = 0.01702127659574468
$ws.Range("E7").Value = 0.00425531914893617
$ws.Range("F7").Value = 0.07659574468085106
$ws.Range("J7").Value = 0.1404255319148936
$ws.Range("O7").Value = 0.02978723404255319
$ws.Range("Q7").Value = 0.1617021276595745
$ws.Range("R7").Value = 0.1106382978723404
$ws.Range("S7").Value = 0.3702127659574468
$ws.Range("B8").Value = 0.09029345372460497
$ws.Range("D8").Value = 0.01580135440180587
$ws.Range("F8").Value = 0.0654627539503386
$ws.Range("J8").Value = 0.1128668171557562
$ws.Range("O8").Value = 0.03837471783295711
$ws.Range("Q8").Value = 0.1647855530474041
$ws.Range("R8").Value = 0.09932279909706546
$ws.Range("S8").Value = 0.4130925507900677
$ws.Range("B9").Value = 0.08888888888888889
$ws.Range("D9").Value = 0.03703703703703703
$ws.Range("F9").Value = 0.03703703703703703
$ws.Range("J9").Value = 0.1037037037037037
$ws.Range("O9").Value = 0.007407407407407408
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("R9").Value = 0.08888888888888889
$ws.Range("S9").Value = 0.4148148148148148
$ws.Range("B10").Value = 0.09808811305070657
$ws.Range("D10").Value = 0.02244389027431421
$ws.Range("E10").Value = 0.0008312551953449709
$ws.Range("F10").Value = 0.06483790523690773
$ws.Range("J10").Value = 0.1014131338320864
$ws.Range("O10").Value = 0.0199501246882793
$ws.Range("Q10").Value = 0.2269326683291771
$ws.Range("R10").Value = 0.1039068994181214
$ws.Range("S10").Value = 0.3615960099750624
$ws.Range("G11").Value = 0.1407185628742515
$ws.Range("J11").Value = 0.0718562874251497
$ws.Range("K11").Value = 0.1736526946107785
$ws.Range("L11").Value = 0.592814371257485
$ws.Range("S11").Value = 0.02095808383233533
$ws.Range("G12").Value = 0.7711442786069652
$ws.Range("J12").Value = 0.1641791044776119
$ws.Range("K12").Value = 0.01990049751243781
$ws.Range("L12").Value = 0.009950248756218905
$ws.Range("S12").Value = 0.03482587064676617
$ws.Range("G13").Value = 0.7017543859649122
$ws.Range("J13").Value = 0.2280701754385965
$ws.Range("S13").Value = 0.07017543859649122
$ws.Range("F15").Value = 0.01904761904761905
$ws.Range("H15").Value = 0.1571428571428571
$ws.Range("I15").Value = 0.04761904761904762
$ws.Range("J15").Value = 0.3571428571428572
$ws.Range("K15").Value = 0.07142857142857142
$ws.Range("M15").Value = 0.01428571428571429
$ws.Range("O15").Value = 0.02857142857142857
$ws.Range("S15").Value = 0.3047619047619048
$ws.Range("F16").Value = 0.00546448087431694
$ws.Range("H16").Value = 0.1693989071038251
$ws.Range("I16").Value = 0.09289617486338798
$ws.Range("J16").Value = 0.3224043715846995
$ws.Range("K16").Value = 0.1475409836065574
$ws.Range("M16").Value = 0.03278688524590164
$ws.Range("O16").Value = 0.08743169398907104
$ws.Range("S16").Value = 0.1420765027322404
$ws.Range("F17").Value = 0.01345291479820628
$ws.Range("H17").Value = 0.1614349775784753
$ws.Range("I17").Value = 0.06053811659192825
$ws.Range("J17").Value = 0.4417040358744395
$ws.Range("K17").Value = 0.1345291479820628
$ws.Range("M17").Value = 0.01345291479820628
$ws.Range("O17").Value = 0.06502242152466367
$ws.Range("S17").Value = 0.1098654708520179
$ws.Range("F18").Value = 0.02643171806167401
$ws.Range("H18").Value = 0.1850220264317181
$ws.Range("I18").Value = 0.07048458149779736
$ws.Range("J18").Value = 0.4008810572687225
$ws.Range("K18").Value = 0.08370044052863436
$ws.Range("M18").Value = 0.03083700440528634
$ws.Range("O18").Value = 0.07048458149779736
$ws.Range("S18").Value = 0.13215859030837
$ws.Range("F19").Value = 0.01269035532994924
$ws.Range("H19").Value = 0.2284263959390863
$ws.Range("I19").Value = 0.05668358714043993
$ws.Range("J19").Value = 0.3764805414551607
$ws.Range("K19").Value = 0.1277495769881557
$ws.Range("M19").Value = 0.02961082910321489
$ws.Range("N19").Value = 0.0008460236886632825
$ws.Range("O19").Value = 0.05668358714043993
$ws.Range("S19").Value = 0.11082910321489
